$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A94").Value = 272
$ws.Range("B94").Value = 172
$ws.Range("C94").Value = 92
$ws.Range("D94").Value = 4
$ws.Range("E94").Value = 4
$ws.Range("F94").Value = 93
$ws.Range("G94").Value = 96
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
